# Weekly update: insert 2 new price records for "Haba" (Femacal de La Calera)
# at rows 73:74 (sorted position by date), pushing the existing rows 73-163
# down to 75-165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 73:74 - this shifts rows 73-163 down to 75-165
# and extends the sheet dimension from A1:R163 to A1:R165 automatically.
$ws.Rows("73:74").Insert()

# --- New row 73 -------------------------------------------------------
$ws.Range("A73").Value = 3
$ws.Range("B73").Value = "Femacal de La Calera"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44763
$ws.Range("E73").Value = 5
$ws.Range("F73").Value = 100112026
$ws.Range("G73").Value = "Haba"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 85
$ws.Range("K73").Value = 16000
$ws.Range("L73").Value = 17000
$ws.Range("M73").Value = 16529
$ws.Range("N73").Value = "$/saco 25 kilos"
$ws.Range("O73").Value = "Provincia de Limarí"
$ws.Range("P73").Value = 661
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"

# --- New row 74 -------------------------------------------------------
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44383
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = 100112026
$ws.Range("G74").Value = "Haba"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 45
$ws.Range("K74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("M74").Value = 13000
$ws.Range("N74").Value = "$/saco 25 kilos"
$ws.Range("O74").Value = "Provincia de Limarí"
$ws.Range("P74").Value = 520
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = "Hortaliza"
